# Trade #97 closed at 2026-02-16 21:39:14 - leadlag UP +0.000%
#
# - Trade #61 (row 50 in "leadlag", row 62 in "All Trades") is closed:
#   exit price / P&L / exit reason / duration are filled in.
# - A brand-new trade #97 is opened and appended as the last row of the
#   "leadlag" sheet.
# - The Summary and Comparison roll-up sheets are refreshed to reflect the
#   new trade counts / win rates / P&L numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Cells.Item(2, 3).Value = 61
$wsSummary.Cells.Item(2, 4).Value = "'67.2%"
$wsSummary.Cells.Item(2, 5).Value = "'+15.9775%"
$wsSummary.Cells.Item(2, 6).Value = "'+0.2619%"
$wsSummary.Cells.Item(3, 5).Value = "'+11.1267%"
$wsSummary.Cells.Item(3, 6).Value = "'+0.1567%"

# ---------------------------------------------------------------------
# leadlag sheet
# ---------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

# Close out trade #61 (row 50): it was OPEN, now CLOSED with a loss.
$wsLeadlag.Cells.Item(50, 7).Value = 68987.194218
$wsLeadlag.Cells.Item(50, 8).Value = "CLOSED"
$wsLeadlag.Cells.Item(50, 9).Value = -0.3256
$wsLeadlag.Cells.Item(50, 10).Value = -3.26
$wsLeadlag.Cells.Item(50, 13).Value = "time_exit_5min"
$wsLeadlag.Cells.Item(50, 14).Value = 5

# Append new trade #97 as row 73 - copy the last existing row first so the
# still-open trade's blank Exit Price / Exit Reason cells and text-typed
# Date column come along untouched, then overwrite just what changed.
$wsLeadlag.Range("A72:N72").Copy()
$wsLeadlag.Range("A73").PasteSpecial()
$wsLeadlag.Cells.Item(73, 1).Value = 97
$wsLeadlag.Cells.Item(73, 3).Value = "21:39:14"
$wsLeadlag.Cells.Item(73, 5).Value = "UP"
$wsLeadlag.Cells.Item(73, 6).Value = 68342.69500000001
$wsLeadlag.Cells.Item(73, 12).Value = "Binance leading with 0.086% move"

# ---------------------------------------------------------------------
# All Trades sheet - mirrors the now-closed trade #61
# ---------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")
$wsAllTrades.Range("A61:N61").Copy()
$wsAllTrades.Range("A62").PasteSpecial()
$wsAllTrades.Cells.Item(62, 1).Value = 61
$wsAllTrades.Cells.Item(62, 3).Value = "21:34:08"
$wsAllTrades.Cells.Item(62, 6).Value = 68763.33
$wsAllTrades.Cells.Item(62, 7).Value = 68987.194218
$wsAllTrades.Cells.Item(62, 9).Value = -0.3256
$wsAllTrades.Cells.Item(62, 10).Value = -3.26
$wsAllTrades.Cells.Item(62, 12).Value = "Coinbase leading with -0.090% move"

# ---------------------------------------------------------------------
# Comparison sheet
# ---------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")
$wsComparison.Cells.Item(2, 4).Value = "2.98"
$wsComparison.Cells.Item(2, 6).Value = "-0.3126%"
